$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.512.76"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "3.039.48"
$ws.Range("E3").Value = "  -2.47%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'209.16"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").Value = "'606.26"
$ws.Range("E6").Value = "  -3.51%  "
$ws.Range("D7").Value = "'0.358"
$ws.Range("E7").Value = "  -9.14%  "
$ws.Range("D8").Value = "'0.882"
$ws.Range("E8").Value = "  +21.81%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "3.039.85"
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("D11").Value = "'0.655"
$ws.Range("E11").Value = "  +18.38%  "
$ws.Range("D12").Value = "'0.186"
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("D13").Value = "'0.0000234"
$ws.Range("E13").Value = "  -6.39%  "
$ws.Range("D14").Value = "'5.32"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").Value = "89.185.46"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.632.61"
$ws.Range("E16").Value = "  -1.58%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'31.68"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").Value = "3.085.25"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").Value = "'3.35"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").Value = "'0.0000203"
$ws.Range("E20").Value = "  -4.46%  "
$ws.Range("D21").Value = "'13.33"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("D22").Value = "'421.49"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").Value = "'4.91"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").Value = "'8.01"
$ws.Range("E24").Value = "  -3.69%  "
$ws.Range("D25").Value = "'5.37"
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("D26").Value = "'83.30"
$ws.Range("E26").Value = "  +5.33%  "
$ws.Range("D27").Value = "'11.51"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("D28").Value = "3.235.51"
$ws.Range("E28").Value = "  -3.47%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +10.15%  "
$ws.Range("D31").Value = "'0.160"
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("D32").Value = "'8.12"
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("D33").Value = "'498.10"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("D34").Value = "'3.53"
$ws.Range("E34").Value = "  -9.01%  "
$ws.Range("D35").Value = "'6.54"
$ws.Range("E35").Value = "  -3.86%  "
$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").Value = "'1.78"
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'22.40"
$ws.Range("E37").Value = "  +2.61%  "
$ws.Range("D38").Value = "'1.22"
$ws.Range("E38").Value = "  -3.17%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "'22.21"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.129"
$ws.Range("E40").Value = "  +3.17%  "
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "'0.138"
$ws.Range("E43").Value = "  +11.25%  "
$ws.Range("D44").Value = "'0.362"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'147.32"
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.79"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'43.31"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0683"
$ws.Range("E48").Value = "  +12.39%  "
$ws.Range("D49").Value = "'4.01"
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("D50").Value = "'1.19"
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("D51").Value = "'154.41"
$ws.Range("E51").Value = "  -6.91%  "
